# Update of "apparel" resource test data:
# swap the sample registrant's name/email from John/JK/john@gmail.com
# to Bharath/Kathir/bharath@gmail.com, refresh the C2 mailto hyperlink's
# display text + tooltip to match (keeping the same mailto: target),
# give the link its "visited" purple/underlined look, and leave the
# active selection on C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the registrant sample values -------------------------------
$ws.Range("A2").Value = "Bharath"
$ws.Range("B2").Value = "Kathir"
$ws.Range("C2").Value = "bharath@gmail.com"

# --- Refresh the hyperlink on C2 ----------------------------------------
# Keep the same mailto: address, but update the display text + tooltip
# so they reflect the new email. Recreate (delete + re-add) rather than
# editing TextToDisplay/ScreenTip in place, since that keeps a single,
# clean hyperlink entry tied to the existing relationship id.
$ws.Range("C2").Hyperlinks.Delete()
$ws.Hyperlinks.Add(
    $ws.Range("C2"),
    "mailto:john@gmail.com",
    [Type]::Missing,
    "mailto:bharath@gmail.com",
    "bharath@gmail.com"
)

# --- Give the link its purple "followed" look ---------------------------
# 0x800080 = RGB(128, 0, 128) packed as (B*65536 + G*256 + R) for the
# COM Font.Color property.
$ws.Range("C2").Font.Color = 8388736
$ws.Range("C2").Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleSingle

# --- Leave the active selection on C2 ------------------------------------
$ws.Range("C2").Select() | Out-Null
